$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 983.2951
$ws.Range("I137").Value = 745.85
$ws.Range("J137").Value = 1435.5714
$ws.Range("K137").Value = 2237.55
$ws.Range("L137").Value = 4306.7142
$ws.Range("M137").Value = 312.4499999999998
$ws.Range("N137").Value = -9406.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7351.506
$ws.Range("I32").Value = 3692.158
$ws.Range("J32").Value = 28744.615
$ws.Range("K32").Value = 3692.158
$ws.Range("L32").Value = 28744.615
$ws.Range("M32").Value = -3405.158
$ws.Range("N32").Value = -29318.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1198.5217
$ws.Range("I74").Value = 1163.0588
$ws.Range("J74").Value = 1299
$ws.Range("K74").Value = 1163.0588
$ws.Range("L74").Value = 1299
$ws.Range("M74").Value = -289.0588
$ws.Range("N74").Value = -3047

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1198.5217
$ws.Range("I77").Value = 1163.0588
$ws.Range("J77").Value = 1299
$ws.Range("K77").Value = 5815.294
$ws.Range("L77").Value = 6495
$ws.Range("M77").Value = -1447.294
$ws.Range("N77").Value = -15231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 45824
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45824
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45824
$ws.Range("N129").Value = -55824

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1516784.8
$ws.Range("I105").Value = 2842421.5
$ws.Range("J105").Value = 1771.2858
$ws.Range("K105").Value = 2842421.5
$ws.Range("L105").Value = 1771.2858
$ws.Range("M105").Value = -2840674.5
$ws.Range("N105").Value = -5265.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6068.3335
$ws.Range("I62").Value = 6068.3335
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 6068.3335
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5444.3335
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 6068.3335
$ws.Range("I65").Value = 6068.3335
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 30341.6675
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -27221.6675
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 893.2564
$ws.Range("I5").Value = 308.06668
$ws.Range("J5").Value = 2843.889
$ws.Range("K5").Value = 924.2000400000001
$ws.Range("L5").Value = 8531.667000000001
$ws.Range("M5").Value = -812.2000400000001
$ws.Range("N5").Value = -8755.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3185.7144
$ws.Range("I63").Value = 2453
$ws.Range("J63").Value = 4162.6665
$ws.Range("K63").Value = 7359
$ws.Range("L63").Value = 12487.9995
$ws.Range("M63").Value = -6610
$ws.Range("N63").Value = -13985.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 3185.7144
$ws.Range("I66").Value = 2453
$ws.Range("J66").Value = 4162.6665
$ws.Range("K66").Value = 22077
$ws.Range("L66").Value = 37463.9985
$ws.Range("M66").Value = -18333
$ws.Range("N66").Value = -44951.9985

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2055.5715
$ws.Range("I70").Value = 1677.8
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 5033.4
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4718.4
$ws.Range("N70").Value = -9630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2055.5715
$ws.Range("I73").Value = 1677.8
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 5033.4
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3941.4
$ws.Range("N73").Value = -11184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1094.4
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1094.4
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 3283.2
$ws.Range("N75").Value = -5279.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1094.4
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1094.4
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 9849.6
$ws.Range("N78").Value = -19833.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 27030652
$ws.Range("I121").Value = 573.3333
$ws.Range("J121").Value = 35718892
$ws.Range("K121").Value = 1719.9999
$ws.Range("L121").Value = 107156676
$ws.Range("M121").Value = -409.9999
$ws.Range("N121").Value = -107159296

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1131.9333
$ws.Range("I129").Value = 534.875
$ws.Range("J129").Value = 1814.2858
$ws.Range("K129").Value = 1604.625
$ws.Range("L129").Value = 5442.857400000001
$ws.Range("M129").Value = 3395.375
$ws.Range("N129").Value = -15442.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2157.4375
$ws.Range("I132").Value = 734.6667
$ws.Range("J132").Value = 2485.7693
$ws.Range("K132").Value = 6612.0003
$ws.Range("L132").Value = 22371.9237
$ws.Range("M132").Value = -4082.0003
$ws.Range("N132").Value = -27431.9237

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 9538.333000000001
$ws.Range("I133").Value = 9930
$ws.Range("J133").Value = 9460
$ws.Range("K133").Value = 29790
$ws.Range("L133").Value = 28380
$ws.Range("M133").Value = -24730
$ws.Range("N133").Value = -38500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3700.0635
$ws.Range("I134").Value = 1290.6666
$ws.Range("J134").Value = 4904.7617
$ws.Range("K134").Value = 3871.9998
$ws.Range("L134").Value = 14714.2851
$ws.Range("M134").Value = 1198.0002
$ws.Range("N134").Value = -24854.2851

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 893.2564
$ws.Range("I135").Value = 308.06668
$ws.Range("J135").Value = 2843.889
$ws.Range("K135").Value = 2772.60012
$ws.Range("L135").Value = 25595.001
$ws.Range("M135").Value = -237.6001200000001
$ws.Range("N135").Value = -30665.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1654.9166
$ws.Range("I140").Value = 928.4666999999999
$ws.Range("J140").Value = 5287.1665
$ws.Range("K140").Value = 2785.4001
$ws.Range("L140").Value = 15861.4995
$ws.Range("M140").Value = 2394.5999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2847.9
$ws.Range("I7").Value = 3928
$ws.Range("J7").Value = 2385
$ws.Range("K7").Value = 3928
$ws.Range("L7").Value = 2385
$ws.Range("M7").Value = -3816
$ws.Range("N7").Value = -2609

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2137.5557
$ws.Range("I40").Value = 1948.0476
$ws.Range("J40").Value = 2800.8333
$ws.Range("K40").Value = 1948.0476
$ws.Range("L40").Value = 2800.8333
$ws.Range("M40").Value = -1812.0476
$ws.Range("N40").Value = -3072.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2712.2
$ws.Range("I122").Value = 2490.2
$ws.Range("J122").Value = 3600.2
$ws.Range("K122").Value = 7470.599999999999
$ws.Range("L122").Value = 10800.6
$ws.Range("M122").Value = -5020.599999999999
$ws.Range("N122").Value = -15700.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2847.9
$ws.Range("I126").Value = 3928
$ws.Range("J126").Value = 2385
$ws.Range("K126").Value = 11784
$ws.Range("L126").Value = 7155
$ws.Range("M126").Value = -9314
$ws.Range("N126").Value = -12095
